# Add manufacturers to parts list (column A, rows 5-25 of Sheet1).
#
# Cells are written in a specific order so that the newly-created shared
# strings land in sharedStrings.xml in the same order as the target
# workbook (Kingbright, TE Connectivity, FTDI (?), N/A, Vishay, ALPS,
# Alpha (Taiwan), Panasonic, KEMET, Murata).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12").Value = "Kingbright"
$ws.Range("A14").Value = "TE Connectivity"
$ws.Range("A23").Value = "FTDI (?)"
$ws.Range("A13").Value = "N/A"
$ws.Range("A7").Value  = "Vishay"
$ws.Range("A17").Value = "ALPS"
$ws.Range("A16").Value = "Alpha (Taiwan)"
$ws.Range("A10").Value = "Panasonic"
$ws.Range("A6").Value  = "KEMET"
$ws.Range("A5").Value  = "Murata"
$ws.Range("A9").Value  = "KEMET"
$ws.Range("A11").Value = "Vishay"
$ws.Range("A15").Value = "TE Connectivity"
$ws.Range("A18").Value = "Kingbright"
$ws.Range("A19").Value = "Vishay"
$ws.Range("A20").Value = "Vishay"
$ws.Range("A21").Value = "TE Connectivity"
$ws.Range("A22").Value = "TE Connectivity"

# A24/A25 are brand-new cells (the row previously started at column B), so
# pull in the grey "sub-part" style (s=7, same as A9-A13) before setting
# their values.
$ws.Range("A9").Copy()
$ws.Range("A24:A25").PasteSpecial(-4122)
$ws.Range("A24").Value = "N/A"
$ws.Range("A25").Value = "N/A"

# Matches the author's final selection/view state in the saved workbook.
$ws.Range("A5").Select()
